$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.04537966666667
$ws.Range("H2").Value = 54.13613900000001
$ws.Range("I2").Value = 0.6797959733292525
$ws.Range("J2").Value = 0.6797959733292525
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.22896466666667
$ws.Range("N2").Value = 93.686894
$ws.Range("O2").Value = 0.2877106972998646
$ws.Range("P2").Value = 0.2877106972998646
$ws.Range("Q2").Value = 563.5385240069185
$ws.Range("R2").Value = 5071.846716062267
$ws.Range("S2").Value = 0.1955845735081994
$ws.Range("T2").Value = 0.1955845735081994

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.04537966666667
$ws.Range("H3").Value = 54.13613900000001
$ws.Range("I3").Value = 0.6797959733292525
$ws.Range("J3").Value = 0.6797959733292525
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.44578266666667
$ws.Range("N3").Value = 121.337348
$ws.Range("O3").Value = 0.3726247238124506
$ws.Range("P3").Value = 0.3726247238124505
$ws.Range("Q3").Value = 729.8595041354858
$ws.Range("R3").Value = 6568.735537219372
$ws.Range("S3").Value = 0.2533087868106287
$ws.Range("T3").Value = 0.2533087868106287

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.04537966666667
$ws.Range("H4").Value = 54.13613900000001
$ws.Range("I4").Value = 0.6797959733292525
$ws.Range("J4").Value = 0.6797959733292525
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.36964133333333
$ws.Range("N4").Value = 76.108924
$ws.Range("O4").Value = 0.2337290805561598
$ws.Range("P4").Value = 0.2337290805561598
$ws.Range("Q4").Value = 457.8048098671596
$ws.Range("R4").Value = 4120.243288804437
$ws.Range("S4").Value = 0.1588880878120259
$ws.Range("T4").Value = 0.1588880878120259

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.04537966666667
$ws.Range("H5").Value = 54.13613900000001
$ws.Range("I5").Value = 0.6797959733292525
$ws.Range("J5").Value = 0.6797959733292525
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.49855033333333
$ws.Range("N5").Value = 34.495651
$ws.Range("O5").Value = 0.1059354983315251
$ws.Range("P5").Value = 0.1059354983315251
$ws.Range("Q5").Value = 207.4957063812766
$ws.Range("R5").Value = 1867.461357431489
$ws.Range("S5").Value = 0.07201452519839853
$ws.Range("T5").Value = 0.07201452519839853

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6001993333333334
$ws.Range("H6").Value = 1.800598
$ws.Range("I6").Value = 0.02261039099934159
$ws.Range("J6").Value = 0.02261039099934159
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.22896466666667
$ws.Range("N6").Value = 93.686894
$ws.Range("O6").Value = 0.2877106972998646
$ws.Range("P6").Value = 0.2877106972998646
$ws.Range("Q6").Value = 18.74360377362356
$ws.Range("R6").Value = 168.692433962612
$ws.Range("S6").Value = 0.006505251360643152
$ws.Range("T6").Value = 0.006505251360643152

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6001993333333334
$ws.Range("H7").Value = 1.800598
$ws.Range("I7").Value = 0.02261039099934159
$ws.Range("J7").Value = 0.02261039099934159
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 40.44578266666667
$ws.Range("N7").Value = 121.337348
$ws.Range("O7").Value = 0.3726247238124506
$ws.Range("P7").Value = 0.3726247238124505
$ws.Range("Q7").Value = 24.27553179267823
$ws.Range("R7").Value = 218.479786134104
$ws.Range("S7").Value = 0.008425190701421181
$ws.Range("T7").Value = 0.008425190701421179

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.6001993333333334
$ws.Range("H8").Value = 1.800598
$ws.Range("I8").Value = 0.02261039099934159
$ws.Range("J8").Value = 0.02261039099934159
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.36964133333333
$ws.Range("N8").Value = 76.108924
$ws.Range("O8").Value = 0.2337290805561598
$ws.Range("P8").Value = 0.2337290805561598
$ws.Range("Q8").Value = 15.22684181517245
$ws.Range("R8").Value = 137.041576336552
$ws.Range("S8").Value = 0.005284705899291381
$ws.Range("T8").Value = 0.005284705899291381

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.6001993333333334
$ws.Range("H9").Value = 1.800598
$ws.Range("I9").Value = 0.02261039099934159
$ws.Range("J9").Value = 0.02261039099934159
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.49855033333333
$ws.Range("N9").Value = 34.495651
$ws.Range("O9").Value = 0.1059354983315251
$ws.Range("P9").Value = 0.1059354983315251
$ws.Range("Q9").Value = 6.901422244366445
$ws.Range("R9").Value = 62.11280019929801
$ws.Range("S9").Value = 0.002395243037985882
$ws.Range("T9").Value = 0.002395243037985882

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.690054333333333
$ws.Range("H10").Value = 14.070163
$ws.Range("I10").Value = 0.1766812397072912
$ws.Range("J10").Value = 0.1766812397072912
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.22896466666667
$ws.Range("N10").Value = 93.686894
$ws.Range("O10").Value = 0.2877106972998646
$ws.Range("P10").Value = 0.2877106972998646
$ws.Range("Q10").Value = 146.4655410604136
$ws.Range("R10").Value = 1318.189869543722
$ws.Range("S10").Value = 0.05083308267598927
$ws.Range("T10").Value = 0.05083308267598927

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.690054333333333
$ws.Range("H11").Value = 14.070163
$ws.Range("I11").Value = 0.1766812397072912
$ws.Range("J11").Value = 0.1766812397072912
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 40.44578266666667
$ws.Range("N11").Value = 121.337348
$ws.Range("O11").Value = 0.3726247238124506
$ws.Range("P11").Value = 0.3726247238124505
$ws.Range("Q11").Value = 189.6929182608582
$ws.Range("R11").Value = 1707.236264347724
$ws.Range("S11").Value = 0.06583579814877075
$ws.Range("T11").Value = 0.06583579814877075

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.690054333333333
$ws.Range("H12").Value = 14.070163
$ws.Range("I12").Value = 0.1766812397072912
$ws.Range("J12").Value = 0.1766812397072912
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 25.36964133333333
$ws.Range("N12").Value = 76.108924
$ws.Range("O12").Value = 0.2337290805561598
$ws.Range("P12").Value = 0.2337290805561598
$ws.Range("Q12").Value = 118.9849962705125
$ws.Range("R12").Value = 1070.864966434612
$ws.Range("S12").Value = 0.04129554370830763
$ws.Range("T12").Value = 0.04129554370830763

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.690054333333333
$ws.Range("H13").Value = 14.070163
$ws.Range("I13").Value = 0.1766812397072912
$ws.Range("J13").Value = 0.1766812397072912
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 11.49855033333333
$ws.Range("N13").Value = 34.495651
$ws.Range("O13").Value = 0.1059354983315251
$ws.Range("P13").Value = 0.1059354983315251
$ws.Range("Q13").Value = 53.92882581790145
$ws.Range("R13").Value = 485.3594323611131
$ws.Range("S13").Value = 0.01871681517422354
$ws.Range("T13").Value = 0.01871681517422354

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.209654333333333
$ws.Range("H14").Value = 9.628962999999999
$ws.Range("I14").Value = 0.1209123959641148
$ws.Range("J14").Value = 0.1209123959641148
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.22896466666667
$ws.Range("N14").Value = 93.686894
$ws.Range("O14").Value = 0.2877106972998646
$ws.Range("P14").Value = 0.2877106972998646
$ws.Range("Q14").Value = 100.2341817678802
$ws.Range("R14").Value = 902.1076359109219
$ws.Range("S14").Value = 0.0347877897550328
$ws.Range("T14").Value = 0.0347877897550328

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.209654333333333
$ws.Range("H15").Value = 9.628962999999999
$ws.Range("I15").Value = 0.1209123959641148
$ws.Range("J15").Value = 0.1209123959641148
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 40.44578266666667
$ws.Range("N15").Value = 121.337348
$ws.Range("O15").Value = 0.3726247238124506
$ws.Range("P15").Value = 0.3726247238124505
$ws.Range("Q15").Value = 129.8169816011249
$ws.Range("R15").Value = 1168.352834410124
$ws.Range("S15").Value = 0.04505494815162994
$ws.Range("T15").Value = 0.04505494815162994

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.209654333333333
$ws.Range("H16").Value = 9.628962999999999
$ws.Range("I16").Value = 0.1209123959641148
$ws.Range("J16").Value = 0.1209123959641148
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 25.36964133333333
$ws.Range("N16").Value = 76.108924
$ws.Range("O16").Value = 0.2337290805561598
$ws.Range("P16").Value = 0.2337290805561598
$ws.Range("Q16").Value = 81.42777924064578
$ws.Range("R16").Value = 732.8500131658119
$ws.Range("S16").Value = 0.02826074313653487
$ws.Range("T16").Value = 0.02826074313653487

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.209654333333333
$ws.Range("H17").Value = 9.628962999999999
$ws.Range("I17").Value = 0.1209123959641148
$ws.Range("J17").Value = 0.1209123959641148
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 11.49855033333333
$ws.Range("N17").Value = 34.495651
$ws.Range("O17").Value = 0.1059354983315251
$ws.Range("P17").Value = 0.1059354983315251
$ws.Range("Q17").Value = 36.90637190443478
$ws.Range("R17").Value = 332.157347139913
$ws.Range("S17").Value = 0.01280891492091719
$ws.Range("T17").Value = 0.01280891492091719
